# edit.ps1 -- Applies the "sentiment analysis integrated with the R data analysis" change.
# Renames/repositions the existing Raw_Data sheet and inserts five new analysis
# sheets (Overall_Sentiment_Correlations, Event_Window_Analysis,
# Correlation_Significance, Daily_Sentiment_Summary, Correlation_Summary) so
# the final tab order is:
#   Volatility_Analysis, Moving_Average_Changes, Statistical_Tests,
#   Overall_Sentiment_Correlations, Event_Window_Analysis,
#   Correlation_Significance, Daily_Sentiment_Summary, Raw_Stock_Data,
#   Correlation_Summary

$wb = $excel.ActiveWorkbook

# Rename Raw_Data -> Raw_Stock_Data (content is untouched).
$rawData = $wb.Worksheets.Item("Raw_Data")
$rawData.Name = "Raw_Stock_Data"

# Insert the first four new sheets right after Statistical_Tests (i.e. before
# Raw_Stock_Data), in order.
$afterSheet = $wb.Worksheets.Item("Statistical_Tests")

# ---- Sheet: Overall_Sentiment_Correlations ----
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Overall_Sentiment_Correlations"
$afterSheet = $ws

$headers = @("Company", "correlation_sentiment", "correlation_vader", "n_observations")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $cell = $ws.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

$rows = @(
    ,@("CNPF", -0.7609799016613027, -0.615471969883548, 9)
    ,@("GSMI", -0.4027456393398635, -0.2822413021159349, 9)
    ,@("JFC", -0.4171755257640853, -0.2097548491947978, 9)
    ,@("MONDE", -0.2487815445421259, -0.1395955539677381, 9)
    ,@("URC", -0.5001911056759699, -0.5322607846706537, 9)
)
$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $v = $row[$c]
        if ($null -ne $v) {
            $cell = $ws.Cells.Item($r, $c + 1)
            if ($v -is [bool]) {
                $cell.Value = $v
            } else {
                $cell.Value = $v
            }
        }
    }
    $r++
}

# ---- Sheet: Event_Window_Analysis ----
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Event_Window_Analysis"
$afterSheet = $ws

$headers = @("Company", "correlation", "avg_return", "avg_sentiment", "n_obs", "event_date")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $cell = $ws.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

$rows = @(
    ,@("CNPF", $null, -4.341736694677882, 0.2775, 1, 45750)
    ,@("GSMI", $null, -1.947623261719669, 0.2775, 1, 45750)
    ,@("JFC", $null, -3.114186851211068, 0.2775, 1, 45750)
    ,@("MONDE", $null, -2.906208718626164, 0.2775, 1, 45750)
    ,@("URC", $null, -1.140413399857465, 0.2775, 1, 45750)
    ,@("CNPF", -0.796052911302524, -0.8175510518010118, -0.009466666666666665, 3, 45848)
    ,@("GSMI", 0.9711599613451505, 0.06910850034554512, -0.009466666666666665, 3, 45848)
    ,@("JFC", -0.6907008212888047, -0.1449163598580867, -0.009466666666666665, 3, 45848)
    ,@("MONDE", 0.9951280168714758, 2.000255770647278, -0.009466666666666665, 3, 45848)
    ,@("URC", -0.6298166871152562, 0.3675537253512813, -0.009466666666666665, 3, 45848)
)
$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $v = $row[$c]
        if ($null -ne $v) {
            $cell = $ws.Cells.Item($r, $c + 1)
            if ($c -eq 5) {
                $cell.Value = $v
                $cell.NumberFormat = "yyyy-mm-dd"
            } elseif ($v -is [bool]) {
                $cell.Value = $v
            } else {
                $cell.Value = $v
            }
        }
    }
    $r++
}

# ---- Sheet: Correlation_Significance ----
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Correlation_Significance"
$afterSheet = $ws

$headers = @("Company", "correlation", "p_value", "significant", "n_observations")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $cell = $ws.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

$rows = @(
    ,@("JFC", -0.4171755257640853, 0.2639428999346807, $false, 9)
    ,@("URC", -0.5001911056759699, 0.1702810694466928, $false, 9)
    ,@("CNPF", -0.7609799016613027, 0.01724163872013969, $true, 9)
    ,@("GSMI", -0.4027456393398635, 0.2824997426247369, $false, 9)
    ,@("MONDE", -0.2487815445421259, 0.5186036260194442, $false, 9)
)
$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $v = $row[$c]
        if ($null -ne $v) {
            $cell = $ws.Cells.Item($r, $c + 1)
            if ($v -is [bool]) {
                $cell.Value = $v
            } else {
                $cell.Value = $v
            }
        }
    }
    $r++
}

# ---- Sheet: Daily_Sentiment_Summary ----
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Daily_Sentiment_Summary"
$afterSheet = $ws

$headers = @("date", "avg_sentiment", "sentiment_count", "avg_vader")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $cell = $ws.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

$rows = @(
    ,@(45631, -0.4398, 1, -0.9883)
    ,@(45702, -0.4232, 1, -0.9136)
    ,@(45743, 0, 1, 0)
    ,@(45750, 0.2775, 2, 0.3488)
    ,@(45753, 0.0682, 1, 0)
    ,@(45794, -0.583, 1, -0.7909)
    ,@(45806, -0.3877, 1, -0.842)
    ,@(45829, 0, 2, 0)
    ,@(45836, -0.2772, 1, -0.2732)
    ,@(45845, -0.1132, 1, -0.2263)
    ,@(45848, 0.1361, 7, 0.2229428571428571)
    ,@(45849, -0.0513, 1, -0.1027)
    ,@(45850, 0, 1, 0)
    ,@(45856, -0.25, 1, 0)
    ,@(45861, 0.07028571428571428, 14, 0.174)
    ,@(45862, 0.541, 1, 0.9712)
)
$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $v = $row[$c]
        if ($null -ne $v) {
            $cell = $ws.Cells.Item($r, $c + 1)
            if ($c -eq 0) {
                $cell.Value = $v
                $cell.NumberFormat = "yyyy-mm-dd"
            } elseif ($v -is [bool]) {
                $cell.Value = $v
            } else {
                $cell.Value = $v
            }
        }
    }
    $r++
}

# Raw_Stock_Data keeps its existing position/content; move insertion point to
# the end of the workbook so Correlation_Summary lands after it.
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---- Sheet: Correlation_Summary ----
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Correlation_Summary"
$afterSheet = $ws

$headers = @("Metric", "Value")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $cell = $ws.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

$rows = @(
    ,@("Companies with |r| > 0.3", "4")
    ,@("Significant correlations (p < 0.05)", "1")
    ,@("Average correlation", "-0.466")
    ,@("Sentiment data points", "16")
    ,@("Date range", "2025-02-14 to 2025-07-23")
)
$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $v = $row[$c]
        if ($null -ne $v) {
            $cell = $ws.Cells.Item($r, $c + 1)
            if ($v -is [bool]) {
                $cell.Value = $v
            } elseif ($c -eq 1) {
                $cell.Value = "'" + $v
            } else {
                $cell.Value = $v
            }
        }
    }
    $r++
}

$wb.Worksheets.Item(1).Activate()
